$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (one decimal point) need to be
# forced to text first, otherwise Excel auto-converts "578.43" etc. to a number.
# Apply a Text number format, set the value, then reset the style back to Normal
# so no stray style survives on the cell (matches original un-styled cells).
$textCells = @("D5", "D6", "D14", "D20", "D21", "D23", "D25", "D26", "D27", "D28", "D31", "D34", "D38", "D39", "D41", "D45", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value updates described by the diff, row by row.
$ws.Range("D2").Value = '62.994.55'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '2.454.83'
$ws.Range("E3").Value = '  +2.24%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '578.43'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").Value = '146.43'
$ws.Range("E6").Value = '  +3.23%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").Value = '2.453.23'
$ws.Range("E9").Value = '  +1.86%  '
$ws.Range("E10").Value = '  +2.97%  '
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("D14").Value = '28.50'
$ws.Range("E14").Value = '  +7.85%  '
$ws.Range("E15").Value = '  +5.59%  '
$ws.Range("D16").Value = '2.898.66'
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").Value = '62.901.74'
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").Value = '2.455.83'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").Value = '11.10'
$ws.Range("E20").Value = '  +4.04%  '
$ws.Range("D21").Value = '331.55'
$ws.Range("E21").Value = '  +2.32%  '
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").Value = '2.09'
$ws.Range("E23").Value = '  +9.29%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").Value = '1.27'
$ws.Range("E25").Value = '  +26.72%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '66.41'
$ws.Range("E26").Value = '  +1.86%  '
$ws.Range("B27").Value = 'Bittensor'
$ws.Range("C27").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D27").Value = '650.73'
$ws.Range("E27").Value = '  +10.97%  '
$ws.Range("D28").Value = '8.58'
$ws.Range("E28").Value = '  +4.44%  '
$ws.Range("E29").Value = '  +6.06%  '
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("E32").Value = '  +6.42%  '
$ws.Range("E33").Value = '  +4.04%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.138'
$ws.Range("E34").Value = '  +4.13%  '
$ws.Range("B35").Value = 'BabyDogeCoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D35").Value = '0.0₆0427'
$ws.Range("E35").Value = '  +52.07%  '
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("D38").Value = '4.79'
$ws.Range("E38").Value = '  +3.61%  '
$ws.Range("D39").Value = '5.55'
$ws.Range("E39").Value = '  +6.02%  '
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").Value = '152.26'
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("E43").Value = '  +10.46%  '
$ws.Range("E44").Value = '  +5.27%  '
$ws.Range("D45").Value = '42.69'
$ws.Range("E45").Value = '  +2.08%  '
$ws.Range("E47").Value = '  +27.51%  '
$ws.Range("D48").Value = '146.48'
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("E49").Value = '  +3.32%  '
$ws.Range("D50").Value = '20.76'
$ws.Range("E50").Value = '  +5.60%  '

# Reset style on the forced-text cells back to Normal (clears the Text number
# format) while keeping the values stored as text.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
